# add user list to project
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# New header cell E1 ("users") - copy the header style (bold/border/centered)
# from the existing "percentage" header (D1) before setting the value.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E1").Value = "users"

# New per-project user lists (column E), aligned with existing rows 2-4
$ws.Range("E2").Value = "['Jonathan Hoff']"
$ws.Range("E3").Value = "['Won Dong Shin']"
$ws.Range("E4").Value = "['Daniel Olivas Hernandez', 'Tanya Verma']"
